$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.931.71'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '2.465.26'
$ws.Range('E3').Value = '  -2.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.46'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.50'
$ws.Range('E6').Value = '  -1.52%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.511'
$ws.Range('D9').Value = '2.465.79'
$ws.Range('E9').Value = '  -2.61%  '
$ws.Range('E10').Value = '  -1.78%  '
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.90'
$ws.Range('E12').Value = '  -2.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.327'
$ws.Range('E13').Value = '  -4.60%  '
$ws.Range('D14').Value = '2.912.98'
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.09'
$ws.Range('E15').Value = '  -3.93%  '
$ws.Range('D16').Value = '66.776.49'
$ws.Range('E16').Value = '  -1.22%  '
$ws.Range('E17').Value = '  -3.96%  '
$ws.Range('D18').Value = '2.458.23'
$ws.Range('E18').Value = '  -1.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.91'
$ws.Range('E19').Value = '  -7.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.35'
$ws.Range('E20').Value = '  -7.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '348.25'
$ws.Range('E21').Value = '  -4.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.00'
$ws.Range('E22').Value = '  -3.61%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.41'
$ws.Range('E24').Value = '  -4.90%  '
$ws.Range('E25').Value = '  -7.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.78'
$ws.Range('E26').Value = '  -4.03%  '
$ws.Range('E27').Value = '  -8.02%  '
$ws.Range('E28').Value = '  -41.68%  '
$ws.Range('E29').Value = '  -2.91%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0891'
$ws.Range('E30').Value = '  -5.48%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '506.82'
$ws.Range('E31').Value = '  -4.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.57'
$ws.Range('E32').Value = '  -8.05%  '
$ws.Range('E33').Value = '  -5.65%  '
$ws.Range('E34').Value = '  -5.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '158.51'
$ws.Range('E36').Value = '  -0.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.114'
$ws.Range('E37').Value = '  -11.12%  '
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.13'
$ws.Range('E39').Value = '  -5.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.32'
$ws.Range('E40').Value = '  -7.85%  '
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('E42').Value = '  -5.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.76'
$ws.Range('E43').Value = '  -5.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.324'
$ws.Range('E44').Value = '  -5.50%  '
$ws.Range('E45').Value = '  -4.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.73'
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '140.64'
$ws.Range('E47').Value = '  -4.95%  '
$ws.Range('E48').Value = '  -7.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.506'
$ws.Range('E49').Value = '  -7.80%  '
$ws.Range('D50').Value = '0.0₆0248'
$ws.Range('E50').Value = '  -9.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0727'
$ws.Range('E51').Value = '  -2.12%  '
